# results diff and heatmaps
#
# - drop the now-unused "donut" sheet (its heatmap run got superseded)
# - tue_logo: remove the stray "x" note in I18 and delete the blank
#   spacer column D so the score columns sit flush against part/util/seam
# - make tue_logo the active sheet/selection (propeller loses focus)

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# 1. Remove the "donut" worksheet entirely.
$donut = $wb.Worksheets.Item("donut")
$donut.Delete()

# 2. Clean up "tue_logo": drop the leftover "x" comment cell, then delete
#    the empty spacer column (column D) so columns E:I shift left to D:H.
$tueLogo = $wb.Worksheets.Item("tue_logo")
$tueLogo.Range("I18").ClearContents()
$tueLogo.Columns.Item(4).Delete()

# 3. tue_logo becomes the active/selected sheet in the saved view.
$tueLogo.Activate()
$tueLogo.Range("I15").Select()
